# Generate Report for Handback
#
# Re-sorts the per-file rows on the "Overview", "zh-cn" and "de-de" sheets
# alphabetically by file GUID, and marks the two files that have now been
# handed back (add9c89b... and c63b682c...) as
#   "Handed back: in sync with en-US"
# filling in their "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the language sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Canonical hyperlink targets, keyed by file guid.
# ---------------------------------------------------------------------
$urls = @{
    "4fcf64be-83aa-4abd-b62d-c144041d2285" = @{
        md        = "https://github.com/OpenLocalizationTest/oltest/blob/51fc3d32cdeadae4350f4ad68c2571f6cadc5789/e2e/4fcf64be-83aa-4abd-b62d-c144041d2285.md"
        zhTarget  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ef4bf02e6bb5b80f58cb07ddec64d3c7cba947d6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4fcf64be-83aa-4abd-b62d-c144041d2285.a96adc780bb03398b9897485a529b36c9626379f.zh-cn.xlf"
        deTarget  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e9dda68fdd57e59be2c6bd2c4713e5f5b9019fa1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4fcf64be-83aa-4abd-b62d-c144041d2285.a96adc780bb03398b9897485a529b36c9626379f.de-de.xlf"
    }
    "6306802d-b0fc-4b02-811f-cf7b96b92804" = @{
        md        = "https://github.com/OpenLocalizationTest/oltest/blob/638416f3782fe79ab8ad333a77fdb8ad905d534d/e2e/6306802d-b0fc-4b02-811f-cf7b96b92804.md"
        zhTarget  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c17495f77f8e924030b11538ff73e9c0567dafbe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/6306802d-b0fc-4b02-811f-cf7b96b92804.247a2058801004946480c8e606d91ee1e060aba9.zh-cn.xlf"
        deTarget  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cf6e07af67271168623784d8ca1d7002de8cef94/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/6306802d-b0fc-4b02-811f-cf7b96b92804.247a2058801004946480c8e606d91ee1e060aba9.de-de.xlf"
    }
    "add9c89b-9271-4bde-bc40-05c6c4229c21" = @{
        md        = "https://github.com/OpenLocalizationTest/oltest/blob/072b352ba31c19be3eb4adfae7dc0bdf2d8a848a/e2e/add9c89b-9271-4bde-bc40-05c6c4229c21.md"
        zhTarget  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c17495f77f8e924030b11538ff73e9c0567dafbe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/add9c89b-9271-4bde-bc40-05c6c4229c21.2846121ea5817640fb35a5985e050e7595de387a.zh-cn.xlf"
        deTarget  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cf6e07af67271168623784d8ca1d7002de8cef94/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/add9c89b-9271-4bde-bc40-05c6c4229c21.2846121ea5817640fb35a5985e050e7595de387a.de-de.xlf"
    }
    "c63b682c-3d16-4e39-a3e2-09927b99e16e" = @{
        md        = "https://github.com/OpenLocalizationTest/oltest/blob/072b352ba31c19be3eb4adfae7dc0bdf2d8a848a/e2e/c63b682c-3d16-4e39-a3e2-09927b99e16e.md"
        zhTarget  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c17495f77f8e924030b11538ff73e9c0567dafbe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/c63b682c-3d16-4e39-a3e2-09927b99e16e.909c8511110cefaec94c046ec90740bd5d624cd4.zh-cn.xlf"
        deTarget  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cf6e07af67271168623784d8ca1d7002de8cef94/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/c63b682c-3d16-4e39-a3e2-09927b99e16e.909c8511110cefaec94c046ec90740bd5d624cd4.de-de.xlf"
    }
}

$HANDED_BACK = "Handed back: in sync with en-US"

# New alphabetical row order (was: 4fcf64be, 6306802d, add9c89b, c63b682c)
$order = @(
    "add9c89b-9271-4bde-bc40-05c6c4229c21",
    "c63b682c-3d16-4e39-a3e2-09927b99e16e",
    "4fcf64be-83aa-4abd-b62d-c144041d2285",
    "6306802d-b0fc-4b02-811f-cf7b96b92804"
)

# ---------------------------------------------------------------------
# 1) "Overview" sheet: File Name / zh-cn / de-de / Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.UsedRange.Hyperlinks.Delete()

$overviewStatus = @{
    "add9c89b-9271-4bde-bc40-05c6c4229c21" = @{ status = $HANDED_BACK;       date = "2016-16-19 10:16:06" }
    "c63b682c-3d16-4e39-a3e2-09927b99e16e" = @{ status = $HANDED_BACK;       date = "2016-16-19 10:16:06" }
    "4fcf64be-83aa-4abd-b62d-c144041d2285" = @{ status = "In Translation";   date = "2016-15-19 10:15:03" }
    "6306802d-b0fc-4b02-811f-cf7b96b92804" = @{ status = "Ready for handoff"; date = "2016-16-19 10:16:06" }
}

$r = 2
foreach ($guid in $order) {
    $fileName = "$guid.md"
    $info = $overviewStatus[$guid]

    $wsOverview.Range("A$r").Value2 = $fileName
    $wsOverview.Range("B$r").Value2 = $info.status
    $wsOverview.Range("C$r").Value2 = $info.status
    $wsOverview.Range("D$r").Value2 = $info.date

    $wsOverview.Hyperlinks.Add($wsOverview.Range("A$r"), $urls[$guid].md, "", "", $fileName) | Out-Null

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Per-language sheets ("zh-cn" / "de-de")
# ---------------------------------------------------------------------
$langSheets = @(
    @{ name = "zh-cn"; targetKey = "zhTarget"; handoffDate = @{
            "add9c89b-9271-4bde-bc40-05c6c4229c21" = "2016-03-19 10:16:02"
            "c63b682c-3d16-4e39-a3e2-09927b99e16e" = "2016-03-19 10:16:02"
            "4fcf64be-83aa-4abd-b62d-c144041d2285" = "2016-03-19 10:14:59"
            "6306802d-b0fc-4b02-811f-cf7b96b92804" = "2016-03-19 10:16:02"
        }
        handbackDate = @{
            "add9c89b-9271-4bde-bc40-05c6c4229c21" = "2016-03-19 10:16:20"
            "c63b682c-3d16-4e39-a3e2-09927b99e16e" = "2016-03-19 10:16:20"
        }
    },
    @{ name = "de-de"; targetKey = "deTarget"; handoffDate = @{
            "add9c89b-9271-4bde-bc40-05c6c4229c21" = "2016-03-19 10:16:06"
            "c63b682c-3d16-4e39-a3e2-09927b99e16e" = "2016-03-19 10:16:06"
            "4fcf64be-83aa-4abd-b62d-c144041d2285" = "2016-03-19 10:15:03"
            "6306802d-b0fc-4b02-811f-cf7b96b92804" = "2016-03-19 10:16:06"
        }
        handbackDate = @{
            "add9c89b-9271-4bde-bc40-05c6c4229c21" = "2016-03-19 10:16:25"
            "c63b682c-3d16-4e39-a3e2-09927b99e16e" = "2016-03-19 10:16:25"
        }
    }
)

$handedBackGuids = @(
    "add9c89b-9271-4bde-bc40-05c6c4229c21",
    "c63b682c-3d16-4e39-a3e2-09927b99e16e"
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.name)
    $ws.UsedRange.Hyperlinks.Delete()

    $r = 2
    foreach ($guid in $order) {
        $fileName = "$guid.md"
        $targetFile = $urls[$guid][$lang.targetKey]
        $isHandedBack = $handedBackGuids -contains $guid

        $status = "Ready for handoff"
        if ($guid -eq "4fcf64be-83aa-4abd-b62d-c144041d2285") { $status = "In Translation" }
        if ($isHandedBack) { $status = $HANDED_BACK }

        $ws.Range("A$r").Value2 = $fileName
        $ws.Range("B$r").Value2 = ".md"
        $ws.Range("C$r").Value2 = $status
        $ws.Range("D$r").Value2 = $targetFile
        $ws.Range("E$r").Value2 = $lang.handoffDate[$guid]

        $ws.Hyperlinks.Add($ws.Range("A$r"), $urls[$guid].md, "", "", $fileName) | Out-Null
        $ws.Hyperlinks.Add($ws.Range("B$r"), $urls[$guid].md, "", "", ".md") | Out-Null
        $ws.Hyperlinks.Add($ws.Range("D$r"), $targetFile, "", "", $targetFile) | Out-Null

        if ($isHandedBack) {
            $ws.Range("F$r").Value2 = $fileName
            $ws.Range("G$r").Value2 = $targetFile
            $ws.Range("H$r").Value2 = $lang.handbackDate[$guid]

            $ws.Hyperlinks.Add($ws.Range("F$r"), $urls[$guid].md, "", "", $fileName) | Out-Null
            $ws.Hyperlinks.Add($ws.Range("G$r"), $targetFile, "", "", $targetFile) | Out-Null
        } else {
            $ws.Range("H$r").Value2 = "0001-01-01 00:00:00"
        }

        $ws.Range("I$r").Value2 = "Include"

        $r = $r + 1
    }
}
